$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.679.97'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.959.35'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '494.88'
$ws.Range('E5').Value = '  -3.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.58'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.28'
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('E10').Value = '  -2.88%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.354'
$ws.Range('E11').Value = '  -1.07%  '
$ws.Range('D12').Value = '3.474.00'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.127'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.63'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '56.833.67'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.02'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').Value = '2.961.14'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.53'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.76'
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '316.36'
$ws.Range('E21').Value = '  -3.70%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '63.02'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E27').Value = '  -6.18%  '
$ws.Range('D28').Value = '0.0₃0882'
$ws.Range('E28').Value = '  -4.20%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.48'
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.01'
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.75'
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -7.14%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.03'
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '154.55'
$ws.Range('E34').Value = '  -1.80%  '
$ws.Range('E35').Value = '  -0.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.69'
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '23.87'
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0661'
$ws.Range('E39').Value = '  -2.84%  '
$ws.Range('D40').Value = '2.994.85'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '37.42'
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.69'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.636'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('D45').Value = '2.192.87'
$ws.Range('E45').Value = '  -4.30%  '
$ws.Range('E46').Value = '  -3.85%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.89'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.932'
$ws.Range('E48').Value = '  -7.49%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0233'
$ws.Range('E49').Value = '  -3.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '19.06'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('E51').Value = '  -10.61%  '
